# "Generate Report for Archive"
#
# Source diff summary:
#   - Shared string "Ready for handoff" -> "In Translation" (status value,
#     used on the Overview sheet in columns E/F (zh-cn / de-de status) and
#     on the per-locale "zh-cn" / "de-de" sheets in column C ("Status")).
#   - Column widths for those same "status" columns shrink to match the
#     shorter replacement text:
#       Overview!E:F   17.2159881591797  -> 13.4101845877511
#       zh-cn!C        17.2159881591797  -> 13.4101845877511
#       de-de!C        17.2159881591797  -> 13.4101845877511

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Each worksheet in this workbook corresponds to a tab named after the
# locale, except the first ("Overview") summary sheet.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: compare with the literal on the LEFT — PowerShell's -eq
            # coerces the right-hand side to the left-hand side's type, and
            # some cells in this sheet hold booleans ("True"/"False"), which
            # would otherwise make a bare `$cell.Value2 -eq $oldStatus`
            # coerce our (non-empty) string to $true and match every
            # True-valued cell too.
            if ($oldStatus -eq [string]$cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Shrink the now-narrower "status" columns to fit the shorter text.
# (ColumnWidth is quantized by the host to 1/6-character steps, so the
# nearest representable width to the target 13.4101845877511 is used.)
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth        # column C (Status)
